$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 22.614608
$ws.Range("H2").Value = 67.843824
$ws.Range("I2").Value = 0.3650188533124966
$ws.Range("J2").Value = 0.3650188533124966
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 92.59233966666666
$ws.Range("N2").Value = 277.777019
$ws.Range("O2").Value = 0.5465415277631132
$ws.Range("P2").Value = 0.5465415277631133
$ws.Range("Q2").Value = 2093.939465364517
$ws.Range("R2").Value = 18845.45518828066
$ws.Range("S2").Value = 0.1994979617517516
$ws.Range("T2").Value = 0.1994979617517516
# Row 3
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 22.614608
$ws.Range("H3").Value = 67.843824
$ws.Range("I3").Value = 0.3650188533124966
$ws.Range("J3").Value = 0.3650188533124966
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 15.14173
$ws.Range("N3").Value = 45.42519
$ws.Range("O3").Value = 0.08937655401050183
$ws.Range("P3").Value = 0.08937655401050183
$ws.Range("Q3").Value = 342.42428839184
$ws.Range("R3").Value = 3081.81859552656
$ws.Range("S3").Value = 0.0326241272579358
$ws.Range("T3").Value = 0.03262412725793579
# Row 4
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 22.614608
$ws.Range("H4").Value = 67.843824
$ws.Range("I4").Value = 0.3650188533124966
$ws.Range("J4").Value = 0.3650188533124966
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 28.17812733333333
$ws.Range("N4").Value = 84.53438199999999
$ws.Range("O4").Value = 0.1663260353686444
$ws.Range("P4").Value = 0.1663260353686444
$ws.Range("Q4").Value = 637.2373038174187
$ws.Range("R4").Value = 5735.135734356767
$ws.Range("S4").Value = 0.06071213870627635
$ws.Range("T4").Value = 0.06071213870627633
# Row 5
$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 22.614608
$ws.Range("H5").Value = 67.843824
$ws.Range("I5").Value = 0.3650188533124966
$ws.Range("J5").Value = 0.3650188533124966
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 33.50281533333334
$ws.Range("N5").Value = 100.508446
$ws.Range("O5").Value = 0.1977558828577406
$ws.Range("P5").Value = 0.1977558828577405
$ws.Range("Q5").Value = 757.6530356597228
$ws.Range("R5").Value = 6818.877320937504
$ws.Range("S5").Value = 0.07218462559653287
$ws.Range("T5").Value = 0.07218462559653285
# Row 6
$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 12.38193366666667
$ws.Range("H6").Value = 37.145801
$ws.Range("I6").Value = 0.1998548561530699
$ws.Range("J6").Value = 0.1998548561530699
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 92.59233966666666
$ws.Range("N6").Value = 277.777019
$ws.Range("O6").Value = 0.5465415277631132
$ws.Range("P6").Value = 0.5465415277631133
$ws.Range("Q6").Value = 1146.472207794135
$ws.Range("R6").Value = 10318.24987014722
$ws.Range("S6").Value = 0.109228978412776
$ws.Range("T6").Value = 0.109228978412776
# Row 7
$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 12.38193366666667
$ws.Range("H7").Value = 37.145801
$ws.Range("I7").Value = 0.1998548561530699
$ws.Range("J7").Value = 0.1998548561530699
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 15.14173
$ws.Range("N7").Value = 45.42519
$ws.Range("O7").Value = 0.08937655401050183
$ws.Range("P7").Value = 0.08937655401050183
$ws.Range("Q7").Value = 187.4838964585767
$ws.Range("R7").Value = 1687.35506812719
$ws.Range("S7").Value = 0.01786233834522593
$ws.Range("T7").Value = 0.01786233834522592
# Row 8
$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 12.38193366666667
$ws.Range("H8").Value = 37.145801
$ws.Range("I8").Value = 0.1998548561530699
$ws.Range("J8").Value = 0.1998548561530699
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 28.17812733333333
$ws.Range("N8").Value = 84.53438199999999
$ws.Range("O8").Value = 0.1663260353686444
$ws.Range("P8").Value = 0.1663260353686444
$ws.Range("Q8").Value = 348.8997034922202
$ws.Range("R8").Value = 3140.097331429982
$ws.Range("S8").Value = 0.03324106587311085
$ws.Range("T8").Value = 0.03324106587311084
# Row 9
$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 12.38193366666667
$ws.Range("H9").Value = 37.145801
$ws.Range("I9").Value = 0.1998548561530699
$ws.Range("J9").Value = 0.1998548561530699
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 33.50281533333334
$ws.Range("N9").Value = 100.508446
$ws.Range("O9").Value = 0.1977558828577406
$ws.Range("P9").Value = 0.1977558828577405
$ws.Range("Q9").Value = 414.8296371039163
$ws.Range("R9").Value = 3733.466733935246
$ws.Range("S9").Value = 0.03952247352195708
$ws.Range("T9").Value = 0.03952247352195707
# Row 10
$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 19.82277733333333
$ws.Range("H10").Value = 59.468332
$ws.Range("I10").Value = 0.3199563508543806
$ws.Range("J10").Value = 0.3199563508543806
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 92.59233966666666
$ws.Range("N10").Value = 277.777019
$ws.Range("O10").Value = 0.5465415277631132
$ws.Range("P10").Value = 0.5465415277631133
$ws.Range("Q10").Value = 1835.4373319847
$ws.Range("R10").Value = 16518.93598786231
$ws.Range("S10").Value = 0.1748694328134638
$ws.Range("T10").Value = 0.1748694328134638
# Row 11
$ws.Range("E11").Value = 3
$ws.Range("G11").Value = 19.82277733333333
$ws.Range("H11").Value = 59.468332
$ws.Range("I11").Value = 0.3199563508543806
$ws.Range("J11").Value = 0.3199563508543806
$ws.Range("K11").Value = 3
$ws.Range("M11").Value = 15.14173
$ws.Range("N11").Value = 45.42519
$ws.Range("O11").Value = 0.08937655401050183
$ws.Range("P11").Value = 0.08937655401050183
$ws.Range("Q11").Value = 300.1511422314533
$ws.Range("R11").Value = 2701.36028008308
$ws.Range("S11").Value = 0.02859659607313962
$ws.Range("T11").Value = 0.02859659607313962
# Row 12
$ws.Range("E12").Value = 3
$ws.Range("G12").Value = 19.82277733333333
$ws.Range("H12").Value = 59.468332
$ws.Range("I12").Value = 0.3199563508543806
$ws.Range("J12").Value = 0.3199563508543806
$ws.Range("K12").Value = 3
$ws.Range("M12").Value = 28.17812733333333
$ws.Range("N12").Value = 84.53438199999999
$ws.Range("O12").Value = 0.1663260353686444
$ws.Range("P12").Value = 0.1663260353686444
$ws.Range("Q12").Value = 558.5687437989803
$ws.Range("R12").Value = 5027.118694190824
$ws.Range("S12").Value = 0.05321707132862811
$ws.Range("T12").Value = 0.05321707132862811
# Row 13
$ws.Range("E13").Value = 3
$ws.Range("G13").Value = 19.82277733333333
$ws.Range("H13").Value = 59.468332
$ws.Range("I13").Value = 0.3199563508543806
$ws.Range("J13").Value = 0.3199563508543806
$ws.Range("K13").Value = 3
$ws.Range("M13").Value = 33.50281533333334
$ws.Range("N13").Value = 100.508446
$ws.Range("O13").Value = 0.1977558828577406
$ws.Range("P13").Value = 0.1977558828577405
$ws.Range("Q13").Value = 664.1188483924525
$ws.Range("R13").Value = 5977.069635532072
$ws.Range("S13").Value = 0.06327325063914903
$ws.Range("T13").Value = 0.06327325063914901
# Row 14
$ws.Range("E14").Value = 3
$ws.Range("G14").Value = 7.135311000000001
$ws.Range("H14").Value = 21.405933
$ws.Range("I14").Value = 0.115169939680053
$ws.Range("J14").Value = 0.115169939680053
$ws.Range("K14").Value = 3
$ws.Range("M14").Value = 92.59233966666666
$ws.Range("N14").Value = 277.777019
$ws.Range("O14").Value = 0.5465415277631132
$ws.Range("P14").Value = 0.5465415277631133
$ws.Range("Q14").Value = 660.6751397393031
$ws.Range("R14").Value = 5946.076257653727
$ws.Range("S14").Value = 0.06294515478512175
$ws.Range("T14").Value = 0.06294515478512175
# Row 15
$ws.Range("E15").Value = 3
$ws.Range("G15").Value = 7.135311000000001
$ws.Range("H15").Value = 21.405933
$ws.Range("I15").Value = 0.115169939680053
$ws.Range("J15").Value = 0.115169939680053
$ws.Range("K15").Value = 3
$ws.Range("M15").Value = 15.14173
$ws.Range("N15").Value = 45.42519
$ws.Range("O15").Value = 0.08937655401050183
$ws.Range("P15").Value = 0.08937655401050183
$ws.Range("Q15").Value = 108.04095262803
$ws.Range("R15").Value = 972.36857365227
$ws.Range("S15").Value = 0.01029349233420049
$ws.Range("T15").Value = 0.01029349233420049
# Row 16
$ws.Range("E16").Value = 3
$ws.Range("G16").Value = 7.135311000000001
$ws.Range("H16").Value = 21.405933
$ws.Range("I16").Value = 0.115169939680053
$ws.Range("J16").Value = 0.115169939680053
$ws.Range("K16").Value = 3
$ws.Range("M16").Value = 28.17812733333333
$ws.Range("N16").Value = 84.53438199999999
$ws.Range("O16").Value = 0.1663260353686444
$ws.Range("P16").Value = 0.1663260353686444
$ws.Range("Q16").Value = 201.059701920934
$ws.Range("R16").Value = 1809.537317288406
$ws.Range("S16").Value = 0.01915575946062914
$ws.Range("T16").Value = 0.01915575946062914
# Row 17
$ws.Range("E17").Value = 3
$ws.Range("G17").Value = 7.135311000000001
$ws.Range("H17").Value = 21.405933
$ws.Range("I17").Value = 0.115169939680053
$ws.Range("J17").Value = 0.115169939680053
$ws.Range("K17").Value = 3
$ws.Range("M17").Value = 33.50281533333334
$ws.Range("N17").Value = 100.508446
$ws.Range("O17").Value = 0.1977558828577406
$ws.Range("P17").Value = 0.1977558828577405
$ws.Range("Q17").Value = 239.053006778902
$ws.Range("R17").Value = 2151.477061010118
$ws.Range("S17").Value = 0.0227755331001016
$ws.Range("T17").Value = 0.0227755331001016
